$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace placeholder text "dooo" with the actual time slot "1*"
$ws.Range("G5").Value = "1*"

# Replace placeholder text "hello" with the actual time slot "11:30-14:30"
$ws.Range("E6").Value = "11:30-14:30"
